$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.389366030693054
$ws.Range("B1").Value = 1.566517472267151
$ws.Range("C1").Value = 1.911339044570923
$ws.Range("D1").Value = 2.608289957046509
$ws.Range("E1").Value = 6.794198989868164
